# Insert a new "Position" column as column B, shifting Age -> C and Nation -> D.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (this shifts existing B:C to C:D and
# copies formatting such as the bold header style from column B into the
# newly inserted column).
$ws.Columns.Item(2).Insert()

# Header for the new column.
$ws.Cells.Item(1, 2).Value = "Position"

$positions = @(
    "Second Striker",
    "Defensive Midfield",
    "Central Midfield",
    "Centre-Back",
    "Central Midfield",
    "Defensive Midfield",
    "Defence",
    "Attacking Midfield",
    "Right Winger",
    "Attacking Midfield",
    "Central Midfield",
    "Left Winger",
    "Centre-Back",
    "Defensive Midfield",
    "Centre-Back",
    "Goalkeeper",
    "Central Midfield",
    "Centre-Forward",
    "Centre-Back",
    "Central Midfield",
    "Attacking Midfield",
    "Centre-Back",
    "Goalkeeper",
    "Central Midfield",
    "Right-Back",
    "Centre-Forward",
    "Centre-Forward",
    "Centre-Forward",
    "Centre-Forward",
    "Left Winger",
    "Second Striker",
    "Defensive Midfield",
    "Defensive Midfield",
    "Sweeper",
    "Defence",
    "Centre-Forward",
    "Defensive Midfield",
    "Defensive Midfield",
    "Defensive Midfield",
    "Attacking Midfield",
    "Second Striker",
    "Central Midfield",
    "Centre-Forward",
    "Right Winger",
    "Centre-Forward",
    "Central Midfield",
    "Attacking Midfield",
    "Centre-Forward",
    "Centre-Back",
    "Goalkeeper",
    "Right Winger",
    "Second Striker",
    "Left Winger",
    "Centre-Back",
    "Centre-Back",
    "Centre-Back",
    "Defensive Midfield",
    "Left Winger",
    "Centre-Back",
    "Left Midfield",
    "Attacking Midfield",
    "Goalkeeper",
    "Left Winger",
    "Right-Back",
    "Left Winger",
    "Midfield",
    "Central Midfield",
    "Attacking Midfield",
    "Centre-Back",
    "Centre-Forward",
    "Attacking Midfield",
    "Midfield",
    "Centre-Back",
    "Centre-Back",
    "Defence"
)

for ($i = 0; $i -lt $positions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $positions[$i]
}
